$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Helper: write a value as literal TEXT (preserving leading zeros /
# avoiding numeric auto-coercion) without leaving a lingering cell
# style behind on the target cell (mirrors cells that only carry
# t="inlineStr"/t="s" with no explicit style index).
# ------------------------------------------------------------------
function Set-TextValue($helper, $range, $text) {
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)
}

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating the existing
#    "总计" sheet (so it inherits sheetPr / styles / page margins),
#    inserted right before "总计" -> ends up between "2021-Q1" and
#    "总计".
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Copy($totalSheet, $null)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"

# NOTE: after Copy(), the *original* $totalSheet variable rebinds to
# the freshly-created copy (COM positional quirk) rather than the
# original "总计" sheet, which has shifted one slot to the right. Get
# a fresh handle onto the real "总计" sheet by name before touching it.
$totalSheet = $wb.Worksheets.Item("总计")

# Helper cell (off to the side, cleared at the end) used to coerce
# numeric-looking strings into real text cells.
$helper1 = $newSheet.Range("ZZ1")
$helper1.NumberFormat = "@"

# Extend the bold/bordered header style (currently only on B1:D1)
# across the new E1:H1 header cells.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Extend the styled index-column (A2) formatting down to the new A3 row.
$newSheet.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)

# Row 2 - 华夏磐锐一年定期开放混合A
$newSheet.Range("A2").Value = 0
Set-TextValue $helper1 $newSheet.Range("B2") "009837"
$newSheet.Range("C2").Value = "华夏磐锐一年定期开放混合A"
Set-TextValue $helper1 $newSheet.Range("D2") "16.45"
Set-TextValue $helper1 $newSheet.Range("E2") "79.44"
Set-TextValue $helper1 $newSheet.Range("F2") "4.60"
Set-TextValue $helper1 $newSheet.Range("G2") "0.7567"
$newSheet.Range("H2").Value = 1

# Row 3 - 华夏磐锐一年定期开放混合C
$newSheet.Range("A3").Value = 1
Set-TextValue $helper1 $newSheet.Range("B3") "009838"
$newSheet.Range("C3").Value = "华夏磐锐一年定期开放混合C"
Set-TextValue $helper1 $newSheet.Range("D3") "0.44"
Set-TextValue $helper1 $newSheet.Range("E3") "79.44"
Set-TextValue $helper1 $newSheet.Range("F3") "4.60"
Set-TextValue $helper1 $newSheet.Range("G3") "0.0202"
$newSheet.Range("H3").Value = 1

$helper1.Clear()

# ------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top data row for 2022-Q1
#    and push the existing 2021-Q1 row down (index column re-numbered).
# ------------------------------------------------------------------
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3:D3").PasteSpecial(-4122)

$helper2 = $totalSheet.Range("ZZ1")
$helper2.NumberFormat = "@"

$totalSheet.Range("A3").Value = 1
Set-TextValue $helper2 $totalSheet.Range("B3") "2021-Q1"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.01

$totalSheet.Range("A2").Value = 0
Set-TextValue $helper2 $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.78

$helper2.Clear()

# Restore the originally active sheet/tab (sheet-copy/rename operations
# above leave the newly created sheet focused).
$wb.Worksheets.Item("2021-Q1").Activate()

Write-Output "Edit complete"
